# Applies the author's edit:
#  1. Slide 1, shape "Google Shape;68;p13" (the "AFT-Net Reconstruction"
#     caption): drop the hyphen so the run reading "AFT-Net Reconstruction"
#     becomes "AFTNet" + " Reconstruction".
#  2. Refresh the auto "datetimeFigureOut" footer field cached text
#     (12/16/2023 -> 10/3/2024) on the slide master and every slide layout.

$p = $ppt.ActivePresentation

# --- 1. "AFT-Net Reconstruction" -> "AFTNet Reconstruction" -------------
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "AFT-Net Reconstruction") {
                # Replace just the "AFT-Net" prefix (7 chars) with "AFTNet"
                # (6 chars, hyphen removed); this naturally splits the run
                # into "AFTNet" + " Reconstruction" while keeping the
                # original run formatting on both pieces.
                $prefix = $sh.TextFrame.TextRange.Characters(1, 7)
                $prefix.Text = "AFTNet"
            }
        }
    }
}

# --- 2. Refresh cached date field text on master + layouts --------------
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.HasTextFrame) {
                if ($sh.TextFrame.TextRange.Text -eq "12/16/2023") {
                    $sh.TextFrame.TextRange.Text = "10/3/2024"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DatePlaceholders $layout.Shapes
}
